# Updated cryptos list on Tue Jan 16 03:35:30 UTC 2024 with GitHub Actions
# Refreshes the coin ranking table (Coin/Link/Price/Volume(1h)) in-place with
# the latest scraped values. Some coins changed rank/position, so their row
# data (Coin name + Link) moved along with their Price/Volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> column letter -> new cell text. Using NumberFormat "@"
# (Text) before assigning values keeps values such as "42.784.82", "1.00" or
# "0.0₃0963" stored verbatim as text instead of being reinterpreted by Excel
# as numbers/dates, matching how the source data is formatted.
$rowUpdates = [ordered]@{
    2 = [ordered]@{ 'D' = '42.784.82'; 'E' = '  +0.41%  ' }
    3 = [ordered]@{ 'D' = '2.522.29'; 'E' = '  +0.25%  ' }
    4 = [ordered]@{ 'D' = '1.00'; 'E' = '  -0.01%  ' }
    5 = [ordered]@{ 'D' = '314.04'; 'E' = '  +1.43%  ' }
    6 = [ordered]@{ 'D' = '95.88'; 'E' = '  -0.68%  ' }
    7 = [ordered]@{ 'E' = '  -1.74%  ' }
    8 = [ordered]@{ 'E' = '  -0.10%  ' }
    9 = [ordered]@{ 'E' = '  -0.98%  ' }
    10 = [ordered]@{ 'E' = '  -1.58%  ' }
    11 = [ordered]@{ 'E' = '  -0.44%  ' }
    12 = [ordered]@{ 'E' = '  -2.64%  ' }
    13 = [ordered]@{ 'E' = '  -3.53%  ' }
    14 = [ordered]@{ 'D' = '2.908.60'; 'E' = '  +0.21%  ' }
    15 = [ordered]@{ 'B' = 'Chainlink'; 'C' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; 'D' = '15.27'; 'E' = '  -3.38%  ' }
    16 = [ordered]@{ 'B' = 'WrappedEther'; 'C' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; 'D' = '2.473.12'; 'E' = '  -1.12%  ' }
    17 = [ordered]@{ 'E' = '  -1.17%  ' }
    18 = [ordered]@{ 'D' = '42.848.08'; 'E' = '  +0.64%  ' }
    19 = [ordered]@{ 'D' = '12.87'; 'E' = '  -0.81%  ' }
    20 = [ordered]@{ 'D' = '6.75'; 'E' = '  +4.42%  ' }
    21 = [ordered]@{ 'D' = '0.0₃0963'; 'E' = '  -1.06%  ' }
    22 = [ordered]@{ 'D' = '69.74'; 'E' = '  -2.54%  ' }
    23 = [ordered]@{ 'D' = '253.73'; 'E' = '  -0.05%  ' }
    24 = [ordered]@{ 'D' = '2.94'; 'E' = '  +0.15%  ' }
    25 = [ordered]@{ 'E' = '  +1.94%  ' }
    26 = [ordered]@{ 'D' = '26.65'; 'E' = '  -1.79%  ' }
    27 = [ordered]@{ 'E' = '  +0.07%  ' }
    28 = [ordered]@{ 'E' = '  +3.50%  ' }
    29 = [ordered]@{ 'D' = '40.94'; 'E' = '  +8.80%  ' }
    30 = [ordered]@{ 'E' = '  +1.81%  ' }
    31 = [ordered]@{ 'E' = '  -0.11%  ' }
    32 = [ordered]@{ 'D' = '157.54'; 'E' = '  +2.29%  ' }
    33 = [ordered]@{ 'D' = '19.54'; 'E' = '  +1.73%  ' }
    34 = [ordered]@{ 'E' = '  +3.43%  ' }
    35 = [ordered]@{ 'E' = '  +2.91%  ' }
    36 = [ordered]@{ 'E' = '  +0.65%  ' }
    37 = [ordered]@{ 'D' = '0.0780'; 'E' = '  -1.11%  ' }
    38 = [ordered]@{ 'D' = '0.111'; 'E' = '  -2.09%  ' }
    39 = [ordered]@{ 'D' = '0.118'; 'E' = '  -1.11%  ' }
    40 = [ordered]@{ 'D' = '23.23'; 'E' = '  -6.79%  ' }
    41 = [ordered]@{ 'E' = '  +13.68%  ' }
    42 = [ordered]@{ 'E' = '  +0.96%  ' }
    43 = [ordered]@{ 'B' = 'RenderToken'; 'C' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; 'D' = '3.80'; 'E' = '  -2.26%  ' }
    44 = [ordered]@{ 'B' = 'NEARProtocol'; 'C' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; 'D' = '3.33'; 'E' = '  -2.11%  ' }
    45 = [ordered]@{ 'E' = '  +0.29%  ' }
    46 = [ordered]@{ 'D' = '2.049.47'; 'E' = '  +0.61%  ' }
    47 = [ordered]@{ 'B' = 'Aave'; 'C' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; 'D' = '109.67'; 'E' = '  +8.00%  ' }
    48 = [ordered]@{ 'B' = 'BitcoinSV'; 'C' = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'; 'D' = '85.48'; 'E' = '  +1.02%  ' }
    49 = [ordered]@{ 'B' = 'FraxShare'; 'C' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; 'D' = '8.92'; 'E' = '  -0.77%  ' }
    50 = [ordered]@{ 'B' = 'ordi'; 'C' = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'; 'D' = '75.75'; 'E' = '  +3.80%  ' }
    51 = [ordered]@{ 'D' = '2.764.40'; 'E' = '  +0.17%  ' }
}

foreach ($row in $rowUpdates.Keys) {
    foreach ($col in $rowUpdates[$row].Keys) {
        $cell = $ws.Range("$col$row")
        $cell.NumberFormat = '@'
        $cell.Value = $rowUpdates[$row][$col]
    }
}
